$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update (GitHub Actions scheduled data refresh)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.704.27"
$ws.Range("E2").Value = "  -7.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.914.04"
$ws.Range("E3").Value = "  -8.17%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.80"
$ws.Range("E5").Value = "  -7.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.50"
$ws.Range("E6").Value = "  -10.75%  "

$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.909.04"
$ws.Range("E8").Value = "  -8.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -4.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.127"
$ws.Range("E10").Value = "  -11.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.85"
$ws.Range("E11").Value = "  -9.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.430"
$ws.Range("E12").Value = "  -5.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000214"
$ws.Range("E13").Value = "  -10.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.52"
$ws.Range("E14").Value = "  -9.93%  "

$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.407.90"
$ws.Range("E16").Value = "  -7.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.948.34"
$ws.Range("E17").Value = "  -7.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "58.823.65"
$ws.Range("E18").Value = "  -7.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.37"
$ws.Range("E19").Value = "  -3.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "416.32"
$ws.Range("E20").Value = "  -9.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").Value = "  -8.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.651"
$ws.Range("E22").Value = "  -6.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.81"
$ws.Range("E23").Value = "  -11.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.54"
$ws.Range("E24").Value = "  -6.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.42"
$ws.Range("E25").Value = "  -7.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.01"
$ws.Range("E27").Value = "  +1.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.46"
$ws.Range("E28").Value = "  -8.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.90"
$ws.Range("E29").Value = "  -9.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.02"
$ws.Range("E30").Value = "  -9.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.97"
$ws.Range("E31").Value = "  -12.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.54"
$ws.Range("E32").Value = "  -9.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0944"
$ws.Range("E33").Value = "  -6.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.41"
$ws.Range("E34").Value = "  -8.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "48.99"
$ws.Range("E35").Value = "  -4.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.900"
$ws.Range("E36").Value = "  -11.94%  "

$ws.Range("E37").Value = "  -19.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.32"
$ws.Range("E38").Value = "  +2.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0621"
$ws.Range("E39").Value = "  -15.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0345"
$ws.Range("E40").Value = "  -11.87%  "

$ws.Range("E41").Value = "  -6.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.622.24"
$ws.Range("E42").Value = "  -6.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "356.32"
$ws.Range("E43").Value = "  -9.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.37"
$ws.Range("E44").Value = "  -10.36%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "118.17"
$ws.Range("E46").Value = "  -6.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.227"
$ws.Range("E47").Value = "  -9.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.106"
$ws.Range("E48").Value = "  -5.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.92"
$ws.Range("E49").Value = "  -9.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.69"
$ws.Range("E50").Value = "  -10.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.95"
$ws.Range("E51").Value = "  -10.42%  "
